$d = $word.ActiveDocument

# wdReplaceAll = 2
$wdReplaceAll = 2

# Merge the run-per-word title into a single run/text node.
$d.Content.Find.Execute(
    "Test 007: Better colours in sepia and night mode of gitbook",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Test 007: Better colours in sepia and night mode of gitbook",
    $wdReplaceAll)

# Merge the run-per-word author line into a single run/text node.
$d.Content.Find.Execute(
    "Emma Cliffe, Skills Centre: MASH, University of Bath",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Emma Cliffe, Skills Centre: MASH, University of Bath",
    $wdReplaceAll)

# Merge the run-per-word date line into a single run/text node.
$d.Content.Find.Execute(
    "September 2020",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "September 2020",
    $wdReplaceAll)
